# Update the first row of the active worksheet ("neg_reaction14") with a new
# set of 46 values (columns A1:AT1), replacing the previous 15 values (A1:O1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(1,2,4,7,8,9,11,15,19,20,23,26,27,29,30,31,32,35,37,38,39,40,43,44,45,46,47,49,50,51,52,53,55,57,58,59,60,61,62,63,65,66,68,69,71,72)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $values[$i]
}
